$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" (column D) values read as plain numbers (or, in one case,
# contain a Unicode subscript digit that the numeric parser mis-handles).
# Force those cells to Text format first so Excel keeps the exact literal
# string instead of silently converting it to a number.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D5").Value = '537.03'
$ws.Range("D6").Value = '145.42'
$ws.Range("D7").Value = '0.998'
$ws.Range("D8").Value = '0.573'
$ws.Range("D9").Value = '7.00'
$ws.Range("D10").Value = '0.102'
$ws.Range("D11").Value = '0.339'
$ws.Range("D15").Value = '21.29'
$ws.Range("D19").Value = '338.35'
$ws.Range("D27").Value = '7.31'
$ws.Range("D28").Value = '0.0₃0755'
$ws.Range("D31").Value = '5.92'
$ws.Range("D32").Value = '18.85'
$ws.Range("D33").Value = '150.87'
$ws.Range("D36").Value = '0.852'
$ws.Range("D37").Value = '0.840'
$ws.Range("D39").Value = '3.62'
$ws.Range("D40").Value = '286.10'
$ws.Range("D44").Value = '0.0539'
$ws.Range("D45").Value = '19.18'
$ws.Range("D46").Value = '0.0943'
$ws.Range("D49").Value = '18.45'
$ws.Range("D50").Value = '4.56'
$ws.Range("D51").Value = '111.34'

# Remaining cells (coin names, links, already-text prices, and all the
# "Volume(1h)" percentage cells) can be set directly.
$ws.Range("D2").Value = '59.465.81'
$ws.Range("E2").Value = '  +0.43%  '
$ws.Range("D3").Value = '2.640.60'
$ws.Range("E3").Value = '  +1.41%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("E5").Value = '  -0.61%  '
$ws.Range("E6").Value = '  +2.84%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("E8").Value = '  +1.26%  '
$ws.Range("E9").Value = '  +8.71%  '
$ws.Range("E10").Value = '  -1.16%  '
$ws.Range("E11").Value = '  +1.02%  '
$ws.Range("D13").Value = '3.109.98'
$ws.Range("E13").Value = '  +1.59%  '
$ws.Range("D14").Value = '59.392.90'
$ws.Range("E14").Value = '  +0.42%  '
$ws.Range("E15").Value = '  +3.45%  '
$ws.Range("D16").Value = '2.670.52'
$ws.Range("E16").Value = '  +1.34%  '
$ws.Range("E17").Value = '  +0.99%  '
$ws.Range("E18").Value = '  +3.05%  '
$ws.Range("E19").Value = '  -0.94%  '
$ws.Range("E20").Value = '  +1.67%  '
$ws.Range("E21").Value = '  -2.27%  '
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("E23").Value = '  -1.88%  '
$ws.Range("E24").Value = '  +2.07%  '
$ws.Range("E25").Value = '  -0.15%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("E27").Value = '  +1.38%  '
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("E30").Value = '  -2.16%  '
$ws.Range("E31").Value = '  +1.48%  '
$ws.Range("E32").Value = '  +0.63%  '
$ws.Range("E33").Value = '  +0.86%  '
$ws.Range("E34").Value = '  +0.55%  '
$ws.Range("E35").Value = '  +2.38%  '
$ws.Range("E36").Value = '  +3.32%  '
$ws.Range("E37").Value = '  +0.56%  '
$ws.Range("E38").Value = '  -1.06%  '
$ws.Range("E39").Value = '  +1.21%  '
$ws.Range("E40").Value = '  +4.21%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("E42").Value = '  +0.74%  '
$ws.Range("E43").Value = '  +0.14%  '
$ws.Range("E44").Value = '  +2.73%  '
$ws.Range("E45").Value = '  +2.69%  '
$ws.Range("E46").Value = '  -1.44%  '
$ws.Range("E47").Value = '  +1.50%  '
$ws.Range("D48").Value = '1.962.51'
$ws.Range("E48").Value = '  +0.41%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("E49").Value = '  -0.59%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("E50").Value = '  +0.98%  '
$ws.Range("E51").Value = '  +0.28%  '
